$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.907.47"
$ws.Range("E2").Value = "  -1.27%  "
$ws.Range("D3").Value = "1.917.90"
$ws.Range("E3").Value = "  +1.16%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "320.07"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -1.44%  "
$ws.Range("E6").Value = "  +0.14%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5054"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -2.35%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4023"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +0.09%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.08316"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -1.15%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "42.59"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -0.30%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.102"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -1.19%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "23.73"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +2.20%  "
$ws.Range("D13").Value = "1.918.47"
$ws.Range("E13").Value = "  +1.95%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.394"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -0.68%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.221"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -1.43%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.004"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +0.27%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "92.17"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -2.11%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001097"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -1.22%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06517"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -1.93%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "18.20"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -0.05%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.000"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +0.02%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.933"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -0.33%  "
$ws.Range("D23").Value = "29.921.44"
$ws.Range("E23").Value = "  -1.15%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.31"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +0.08%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.195"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -1.52%  "
$ws.Range("B26").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C26").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D26").Value = "2.137.32"
$ws.Range("E26").Value = "  +1.46%  "
$ws.Range("B27").Value = "EthereumClassic"
$ws.Range("C27").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "22.08"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +1.99%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "162.10"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +0.18%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.315"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -1.66%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "128.92"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -0.32%  "
$ws.Range("E31").Value = "  +3.56%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.1036"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -1.78%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.961"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -2.20%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.805"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +1.52%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.02448"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -1.88%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.402"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +2.13%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06412"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -2.14%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2152"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -2.01%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.6490"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -0.14%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.697"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -0.84%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.192"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -2.31%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "11.35"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -3.84%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.214"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -1.18%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.226"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +8.29%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.34"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +1.03%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.6052"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -0.63%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.635"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -1.43%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.208"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -2.29%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "121.94"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -2.12%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "78.86"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -0.32%  "
$ws.Range("E51").Value = "  -2.52%  "
